# Reverse the order of the comma-separated "Recorded By" entries in column G
# for every data row (row 1 is the header) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2
    if ($val -ne $null -and $val.Contains(",")) {
        $parts = $val -split ", "

        # Build the reversed list manually (array helper methods are unreliable here)
        $reversed = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }

        $newVal = [string]::Join(", ", $reversed)
        $cell.Value2 = $newVal
    }
}
